$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Main Data")

# Row 17: new flight entry for Friday, Jan 13
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Friday, Jan 13"
$ws.Range("C17").Value = "6:05 PM"
$ws.Range("D17").Value = "FR2469"
$ws.Range("E17").Value = "London"
$ws.Range("F17").Value = "(STN)"
$ws.Range("G17").Value = "Lauda Europe "
$ws.Range("H17").Value = "A320"
$ws.Range("I17").Value = "(9H-LOA)"
$ws.Range("J17").Value = "7:51 PM"
$ws.Range("K17").Value = ""
$ws.Range("L17").Value = "1 hours, 46 minutes"
$ws.Range("M17").Value = ""

# Row 18: new flight entry for Friday, Jan 13
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Friday, Jan 13"
$ws.Range("C18").Value = "7:25 PM"
$ws.Range("D18").Value = "FR1979"
$ws.Range("E18").Value = "Dublin"
$ws.Range("F18").Value = "(DUB)"
$ws.Range("G18").Value = "Ryanair "
$ws.Range("H18").Value = "B738"
$ws.Range("I18").Value = "(EI-DHZ)"
$ws.Range("J18").Value = "7:39 PM"
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = "0 hours, 14 minutes"
$ws.Range("M18").Value = ""
